# Automatic update of files.
# Rows 10-21 on the "Artfynd" sheet get their observation records
# re-shuffled into a new order (rows keep the same set of underlying
# records, but each row now shows the record that used to live in a
# different row). Additionally:
#   - Ost/Nord (Q/R) coordinates are stored rounded to whole metres
#     instead of the original full-precision decimal values.
#   - Starttid/Sluttid (Z/AB) are cleared out entirely.
#   - Publik kommentar (AC) is only populated ("ringhack äldre") for
#     the three "Tretåig hackspett" (Picoides tridactylus) records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target content for rows 10..21, taken from the record that used to be
# in the row named by "Src" (rounding Q/R to the nearest integer).
# "Bird" marks the three "Tretåig hackspett" (Picoides tridactylus)
# records, which are the only ones carrying (empty)
# Ålder-Stadium/Kön/Aktivitet/Metod (K/L/M/N) cells and a Publik
# kommentar (AC) of "ringhack äldre".
$rows = @(
    @{ Row = 10; A = 111936796; B = 56398;  D = "NT"; E = 100109; F = "Tretåig hackspett"; G = "Picoides tridactylus";      H = "(Linnaeus, 1758)";            Q = 448883; R = 7087229; Bird = $true },
    @{ Row = 11; A = 111936858; B = 89845;  D = "VU"; E = 1209;   F = "Rynkskinn";          G = "Phlebia centrifuga";        H = "P.Karst.";                    Q = 448737; R = 7087496; Bird = $false },
    @{ Row = 12; A = 111936866; B = 89423;  D = "NT"; E = 5432;   F = "Granticka";          G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 448766; R = 7087417; Bird = $false },
    @{ Row = 13; A = 111936870; B = 89423;  D = "NT"; E = 5432;   F = "Granticka";          G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 449019; R = 7087277; Bird = $false },
    @{ Row = 14; A = 111936798; B = 56398;  D = "NT"; E = 100109; F = "Tretåig hackspett"; G = "Picoides tridactylus";      H = "(Linnaeus, 1758)";            Q = 448923; R = 7087371; Bird = $true },
    @{ Row = 15; A = 111936792; B = 90087;  D = "LC"; E = 3298;   F = "Trådticka";          G = "Climacocystis borealis";    H = "(Fr.) Kotl. & Pouzar";         Q = 448761; R = 7087579; Bird = $false },
    @{ Row = 16; A = 111936868; B = 89423;  D = "NT"; E = 5432;   F = "Granticka";          G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 448988; R = 7087187; Bird = $false },
    @{ Row = 17; A = 111936867; B = 89423;  D = "NT"; E = 5432;   F = "Granticka";          G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 448792; R = 7087386; Bird = $false },
    @{ Row = 18; A = 111936893; B = 77515;  D = "NT"; E = 6425;   F = "Garnlav";            G = "Alectoria sarmentosa";      H = "(Ach.) Ach.";                  Q = 448742; R = 7087502; Bird = $false },
    @{ Row = 19; A = 111936795; B = 56398;  D = "NT"; E = 100109; F = "Tretåig hackspett"; G = "Picoides tridactylus";      H = "(Linnaeus, 1758)";            Q = 448749; R = 7087422; Bird = $true },
    @{ Row = 20; A = 111936865; B = 89423;  D = "NT"; E = 5432;   F = "Granticka";          G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 448738; R = 7087426; Bird = $false },
    @{ Row = 21; A = 111936869; B = 89423;  D = "NT"; E = 5432;   F = "Granticka";          G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 449144; R = 7087118; Bird = $false }
)

foreach ($r in $rows) {
    $n = $r.Row

    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R

    # Starttid / Sluttid are removed entirely on every affected row.
    $ws.Range("Z$n").ClearContents()
    $ws.Range("AB$n").ClearContents()

    if ($r.Bird) {
        # Bird records carry empty (but present) Ålder-Stadium / Kön /
        # Aktivitet / Metod cells plus a Publik kommentar. A leading
        # apostrophe forces Excel to materialize a real, empty text
        # cell instead of leaving it blank/absent.
        $ws.Range("K$n").Value = "'"
        $ws.Range("L$n").Value = "'"
        $ws.Range("M$n").Value = "'"
        $ws.Range("N$n").Value = "'"
        $ws.Range("AC$n").Value = "ringhack äldre"
    } else {
        $ws.Range("K$n").ClearContents()
        $ws.Range("L$n").ClearContents()
        $ws.Range("M$n").ClearContents()
        $ws.Range("N$n").ClearContents()
        $ws.Range("AC$n").ClearContents()
    }
}
